# Regenerate orders with updated distance/size codes.
# The experiment's distance conditions and one size condition were renamed:
#   D80 -> D86
#   D51 -> D55
#   D64 -> D69
#   S30 -> S31
# These tokens appear embedded inside many cell values across the sheet
# (Condition, Filename_Left, Filename_Right, Distance, Size, ...), so we do
# an in-place substring find & replace over the whole used range for each
# token. None of the replacement targets collide with any other token's
# source text, so the four replacements are order independent and safe to
# run as simple sequential substring substitutions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.UsedRange

$rng.Replace("D80", "D86") | Out-Null
$rng.Replace("D51", "D55") | Out-Null
$rng.Replace("D64", "D69") | Out-Null
$rng.Replace("S30", "S31") | Out-Null
